$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 240 (old rows 240-263 shift down
# to 242-265). The new rows 240-241 carry a brand-new weekly price record
# ("$/bandeja 6 kilos") while the previously existing records are preserved,
# just displaced two rows down.
$ws.Range("A240:A241").EntireRow.Insert()

# --- New row 240: Primera, $/bandeja 6 kilos ---
$ws.Range("A240").Value2 = 6
$ws.Range("B240").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C240").Value2 = "Metropolitana"
$ws.Range("D240").Value2 = 44826
$ws.Range("E240").Value2 = 13
$ws.Range("F240").Value2 = "Fruta"
$ws.Range("G240").Value2 = 100107
$ws.Range("H240").Value2 = "Otros"
$ws.Range("I240").Value2 = 100107002
$ws.Range("J240").Value2 = "Chirimoya"
$ws.Range("K240").Value2 = "Cultivar IV Región"
$ws.Range("L240").Value2 = "Primera"
$ws.Range("M240").Value2 = 275
$ws.Range("N240").Value2 = 26000
$ws.Range("O240").Value2 = 26000
$ws.Range("P240").Value2 = 26000
$ws.Range("Q240").Value2 = "$/bandeja 6 kilos"
$ws.Range("R240").Value2 = "Provincia de Limarí"
$ws.Range("S240").Value2 = 4333
$ws.Range("T240").Value2 = 6

# --- New row 241: Segunda, $/bandeja 6 kilos ---
$ws.Range("A241").Value2 = 6
$ws.Range("B241").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C241").Value2 = "Metropolitana"
$ws.Range("D241").Value2 = 44826
$ws.Range("E241").Value2 = 13
$ws.Range("F241").Value2 = "Fruta"
$ws.Range("G241").Value2 = 100107
$ws.Range("H241").Value2 = "Otros"
$ws.Range("I241").Value2 = 100107002
$ws.Range("J241").Value2 = "Chirimoya"
$ws.Range("K241").Value2 = "Cultivar IV Región"
$ws.Range("L241").Value2 = "Segunda"
$ws.Range("M241").Value2 = 275
$ws.Range("N241").Value2 = 22000
$ws.Range("O241").Value2 = 22000
$ws.Range("P241").Value2 = 22000
$ws.Range("Q241").Value2 = "$/bandeja 6 kilos"
$ws.Range("R241").Value2 = "Provincia de Limarí"
$ws.Range("S241").Value2 = 3667
$ws.Range("T241").Value2 = 6
